# MAI_holdings.xlsx daily refresh: bump the "as of" date in the disclosure
# banner and update the Weight / Percent Change figures for the current
# holdings. The sheet ships protected, so it must be unprotected before the
# cells can be written and protected again afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Disclosure banner: "as of" date rolls from 2021-05-10 -> 2021-05-11
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."
# Re-fit the row so the wrapped banner text doesn't leave a stale custom height behind
$ws.Rows.Item(10).AutoFit()

# Holdings table: refreshed Weight (D) / Percent Change (E) figures
$ws.Range("D2").Value = 0.4765519535147614
$ws.Range("E2").Value = -0.00663026521060861

$ws.Range("D3").Value = 0.3417990530784732
$ws.Range("E3").Value = -0.01125476699841887

$ws.Range("D4").Value = 0.09648230391485531
$ws.Range("E4").Value = -0.01312987934705467

$ws.Range("D5").Value = 0.05327044978772523
$ws.Range("E5").Value = -0.001032939286124024

$ws.Range("D6").Value = 0.03189623970418484
$ws.Range("E6").Value = -0.01236342725704431

$ws.Range("E7").Value = -0.008722707530373963

$ws.Protect()
